# Applies updated symbol-list values to the cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = '''292.53'
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = '''-7.05%'
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = '''40.65'
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = '''-0.88%'
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.Value = '''5.043'
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = '''-1.55%'
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = '''0.07401'
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = '''-3.04%'
$cell.Style = "Normal"
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cell = $ws.Range("D6")
$cell.Value = '''1.551'
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = '''-8.10%'
$cell.Style = "Normal"
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D7")
$cell.Value = '''0.9245'
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = '''-1.21%'
$cell.Style = "Normal"
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$cell = $ws.Range("D8")
$cell.Value = '''2.349'
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = '''-3.13%'
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = '''0.1152'
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = '''-8.22%'
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = '''0.1726'
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = '''-5.62%'
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = '''0.08668'
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = '''-4.10%'
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.Value = '''0.04174'
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = '''1.01%'
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = '''0.1056'
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = '''0.06%'
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = '''0.001245'
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = '''-1.96%'
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = '''0.005969'
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = '''2.69%'
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = '''3.414'
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = '''1.49%'
$cell.Style = "Normal"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$cell = $ws.Range("D17")
$cell.Value = '''4.280'
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = '''-1.24%'
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = '''-2.27%'
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.Value = '''7.658'
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = '''-9.22%'
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = '''2.20%'
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = '''4.85%'
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = '''0.03857'
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = '''-4.53%'
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = '''0.001257'
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = '''-0.62%'
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = '''0.003866'
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = '''-4.55%'
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = '''0.0001279'
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = '''0.31%'
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.Value = '''0.02345'
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = '''-5.63%'
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = '''0.05003'
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = '''-3.75%'
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = '''169.77%'
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = '''0.007695'
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = '''-1.17%'
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = '''0.1287'
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = '''-1.06%'
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = '''0.007321'
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = '''-0.67%'
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = '''0.007109'
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = '''-12.92%'
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = '''0.3152'
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = '''0.32%'
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = '''0.00006409'
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = '''-3.71%'
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = '''0.00000000749'
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = '''-0.47%'
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = '''0.01706'
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = '''-94.15%'
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.Value = '''0.00002098'
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = '''-0.47%'
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = '''0.0001998'
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = '''-0.47%'
$cell.Style = "Normal"
